$wb = $excel.ActiveWorkbook

$wsProduto = $wb.Worksheets.Item("Product Backlog")
$wsSprint  = $wb.Worksheets.Item("Sprint Backlog")

# --- Content changes -------------------------------------------------

# "Criar relatório de clientes cadastrados" (item 16) was reworked into a
# new requirement: a CNH-based customer search.
$wsProduto.Range("D21").Value = "Alterado - Criar pesquisa de cliente por CNH"
$wsSprint.Range("C23").Value  = "Alterado - Criar pesquisa de cliente por CNH"

# Realization dates filled in for tasks that were completed (numeric date
# serials keep the existing date-formatted cell style untouched).
$wsProduto.Range("G12").Value = 42517   # 27-May-16
$wsProduto.Range("G13").Value = 42516   # 26-May-16
$wsProduto.Range("G21").Value = 42517   # 27-May-16

# Corresponding "Feito" (done) marks on the Sprint Backlog sheet.
$wsSprint.Range("E14").Value = "X"
$wsSprint.Range("E15").Value = "X"
$wsSprint.Range("E23").Value = "X"

# --- Active sheet / selection -----------------------------------------

# Set the (soon to be inactive) Sprint Backlog's selection first --
# selecting a range activates its sheet, so this must happen before we
# activate + select on the Product Backlog sheet below.
$wsSprint.Range("E21").Select() | Out-Null

# Product Backlog becomes the active/selected sheet with D15 selected.
$wsProduto.Activate() | Out-Null
$wsProduto.Range("D15").Select() | Out-Null
